# Edit: (1) change the table style on the financial-documents table (slide 5)
# to the built-in style {CC20FF9D-34F6-4F28-9EF9-BE0484AF79EB}; (2) swap the
# deck's active color scheme from the colourful "Integral"/"Red Violet"
# palette to the plain "Office Theme" palette (this is the practical,
# visible effect of the theme1.xml / theme2.xml swap recorded upstream).

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------------
$s = $p.Slides.Item(5)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{CC20FF9D-34F6-4F28-9EF9-BE0484AF79EB}")
    }
}

# --- 2) Theme colours -------------------------------------------------------
# msoColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$cs = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $bgr = $officeThemeColors[$i - 1]
    $r = ($bgr -band 0xFF0000) / 0x10000
    $g = ($bgr -band 0x00FF00) / 0x100
    $b = ($bgr -band 0x0000FF)
    $rgb = $r + ($g * 256) + ($b * 65536)
    $cs.Colors($i).RGB = $rgb
}
